$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 144
$ws.Cells.Item(144, 1).Value = 6
$ws.Cells.Item(144, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(144, 3).Value = "Metropolitana"
$ws.Cells.Item(144, 4).Value = 44890
$ws.Cells.Item(144, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(144, 5).Value = 13
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100103
$ws.Cells.Item(144, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(144, 9).Value = 100103003
$ws.Cells.Item(144, 10).Value = "Damasco"
$ws.Cells.Item(144, 11).Value = "Castle Brite"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 250
$ws.Cells.Item(144, 14).Value = 21000
$ws.Cells.Item(144, 15).Value = 22000
$ws.Cells.Item(144, 16).Value = 21500
$ws.Cells.Item(144, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(144, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 19).Value = 1344
$ws.Cells.Item(144, 20).Value = 16

# Row 145
$ws.Cells.Item(145, 1).Value = 6
$ws.Cells.Item(145, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(145, 3).Value = "Metropolitana"
$ws.Cells.Item(145, 4).Value = 44890
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(145, 5).Value = 13
$ws.Cells.Item(145, 6).Value = "Fruta"
$ws.Cells.Item(145, 7).Value = 100103
$ws.Cells.Item(145, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(145, 9).Value = 100103003
$ws.Cells.Item(145, 10).Value = "Damasco"
$ws.Cells.Item(145, 11).Value = "Castle Brite"
$ws.Cells.Item(145, 12).Value = "Segunda"
$ws.Cells.Item(145, 13).Value = 185
$ws.Cells.Item(145, 14).Value = 18000
$ws.Cells.Item(145, 15).Value = 18000
$ws.Cells.Item(145, 16).Value = 18000
$ws.Cells.Item(145, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(145, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(145, 19).Value = 1125
$ws.Cells.Item(145, 20).Value = 16
